$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    ,@('D2', '36.499.26', $false)
    ,@('E2', '  +0.22%  ', $false)
    ,@('D3', '1.953.84', $false)
    ,@('E3', '  +0.50%  ', $false)
    ,@('E4', '  -0.02%  ', $false)
    ,@('D5', '243.96', $true)
    ,@('E5', '  +0.53%  ', $false)
    ,@('D6', '0.613', $true)
    ,@('E6', '  +0.21%  ', $false)
    ,@('D7', '58.08', $true)
    ,@('E7', '  +0.99%  ', $false)
    ,@('E8', '  -0.01%  ', $false)
    ,@('D9', '0.376', $true)
    ,@('E9', '  +3.68%  ', $false)
    ,@('D10', '0.0787', $true)
    ,@('E10', '  -8.15%  ', $false)
    ,@('E11', '  -1.13%  ', $false)
    ,@('D12', '14.03', $true)
    ,@('E12', '  +3.77%  ', $false)
    ,@('D13', '0.836', $true)
    ,@('E13', '  +2.28%  ', $false)
    ,@('D14', '2.242.61', $false)
    ,@('E14', '  +0.61%  ', $false)
    ,@('D15', '21.15', $true)
    ,@('E15', '  -0.64%  ', $false)
    ,@('D16', '5.28', $true)
    ,@('E16', '  +1.78%  ', $false)
    ,@('D17', '1.958.20', $false)
    ,@('E17', '  +1.08%  ', $false)
    ,@('D18', '36.423.98', $false)
    ,@('E18', '  +0.19%  ', $false)
    ,@('D19', '69.70', $true)
    ,@('E19', '  +0.48%  ', $false)
    ,@('D20', '0.0₃0846', $false)
    ,@('E20', '  -3.21%  ', $false)
    ,@('D21', '229.25', $true)
    ,@('E21', '  +0.23%  ', $false)
    ,@('D22', '5.03', $true)
    ,@('E22', '  +0.32%  ', $false)
    ,@('D23', '0.999', $true)
    ,@('E23', '  -0.12%  ', $false)
    ,@('D24', '2.43', $true)
    ,@('E24', '  +1.15%  ', $false)
    ,@('E25', '  +2.49%  ', $false)
    ,@('D26', '9.11', $true)
    ,@('E26', '  -1.30%  ', $false)
    ,@('D27', '0.138', $true)
    ,@('E27', '  +3.07%  ', $false)
    ,@('D28', '160.27', $true)
    ,@('E28', '  -0.52%  ', $false)
    ,@('D29', '19.33', $true)
    ,@('E29', '  +0.15%  ', $false)
    ,@('D30', '0.120', $true)
    ,@('E30', '  +1.61%  ', $false)
    ,@('E31', '  +3.88%  ', $false)
    ,@('D32', '4.73', $true)
    ,@('E32', '  +2.24%  ', $false)
    ,@('D33', '0.0609', $true)
    ,@('E33', '  -4.10%  ', $false)
    ,@('D34', '4.40', $true)
    ,@('E34', '  +3.72%  ', $false)
    ,@('E35', '  +12.77%  ', $false)
    ,@('E36', '  -0.06%  ', $false)
    ,@('E37', '  +5.11%  ', $false)
    ,@('E38', '  -1.59%  ', $false)
    ,@('D39', '5.21', $true)
    ,@('D40', '0.0976', $true)
    ,@('E40', '  +0.08%  ', $false)
    ,@('E41', '  +1.85%  ', $false)
    ,@('E42', '  -0.37%  ', $false)
    ,@('D43', '0.0209', $true)
    ,@('E43', '  -0.63%  ', $false)
    ,@('B44', 'Maker', $false)
    ,@('C44', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', $false)
    ,@('D44', '1.368.71', $false)
    ,@('E44', '  +1.75%  ', $false)
    ,@('B45', 'InjectiveProtocol', $false)
    ,@('C45', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj', $false)
    ,@('D45', '15.76', $true)
    ,@('E45', '  +0.15%  ', $false)
    ,@('B46', 'Aave', $false)
    ,@('C46', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', $false)
    ,@('D46', '87.68', $true)
    ,@('E46', '  +0.37%  ', $false)
    ,@('B47', 'ARBITRUM', $false)
    ,@('C47', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', $false)
    ,@('D47', '1.02', $true)
    ,@('E47', '  -0.30%  ', $false)
    ,@('D48', '7.16', $true)
    ,@('E48', '  +0.09%  ', $false)
    ,@('D49', '2.83', $true)
    ,@('E49', '  +0.04%  ', $false)
    ,@('D50', '2.132.91', $false)
    ,@('E50', '  +0.63%  ', $false)
    ,@('D51', '43.87', $true)
    ,@('E51', '  -1.23%  ', $false)
)

foreach ($change in $changes) {
    $cellRef = $change[0]
    $val = $change[1]
    $isNumericLooking = $change[2]
    $rng = $ws.Range($cellRef)
    if ($isNumericLooking) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $val
}